$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D, shifting D:K to E:L
$ws.Range("D1").EntireColumn.Insert()

# Copy number formats from column E (original D, now shifted) into new column D
$ws.Range("E1:E102").Copy()
$ws.Range("D1:D102").PasteSpecial(-4122)

# Header date row (row 7) now duplicated into rows 7, 38, 80 - set new column D date
$ws.Range("D7").Value = 43404
$ws.Range("D38").Value = 43404
$ws.Range("D80").Value = 43404

# New column D quarter data for each section
$ws.Range("D8").Value = 116500
$ws.Range("D9").Value = 86900
$ws.Range("D10").Value = 29600
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 98100
$ws.Range("D18").Value = 18400
$ws.Range("D20").Value = 1400
$ws.Range("D21").Value = 21000
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 19800
$ws.Range("D24").Value = -12600
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 32400
$ws.Range("D27").Value = 32400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -1400
$ws.Range("D33").Value = 32400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 32400

$ws.Range("D41").Value = 155800
$ws.Range("D42").Value = 159000
$ws.Range("D43").Value = 99200
$ws.Range("D44").Value = "NA"
$ws.Range("D45").Value = 25500
$ws.Range("D46").Value = 439500
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 19900
$ws.Range("D49").Value = 40700
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 700
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 500800
$ws.Range("D57").Value = 65000
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 34900
$ws.Range("D60").Value = 99900
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 1400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 101200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 253700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 399600
$ws.Range("D77").Value = 0

$ws.Range("D81").Value = 32400
$ws.Range("D83").Value = 1200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -40600
$ws.Range("D91").Value = -2000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 34500
$ws.Range("D96").Value = -3900
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -3900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -10000
